$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 772.36365
$ws.Range("I2").Value = 312.375
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 312.375
$ws.Range("L2").Value = 1999
$ws.Range("M2").Value = -199.375
$ws.Range("N2").Value = -2225
$ws.Range("H4").Value = 1300.7273
$ws.Range("I4").Value = 425.75
$ws.Range("K4").Value = 425.75
$ws.Range("M4").Value = -311.75
$ws.Range("H21").Value = 26599.6
$ws.Range("I21").Value = 26599.6
$ws.Range("K21").Value = 26599.6
$ws.Range("M21").Value = -26131.6
$ws.Range("H23").Value = 26599.6
$ws.Range("I23").Value = 26599.6
$ws.Range("K23").Value = 26599.6
$ws.Range("M23").Value = -26365.6
$ws.Range("H33").Value = 367.33334
$ws.Range("I33").Value = 264.75
$ws.Range("J33").Value = 695.6
$ws.Range("K33").Value = 264.75
$ws.Range("L33").Value = 695.6
$ws.Range("M33").Value = -35.75
$ws.Range("N33").Value = -1153.6
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H51").Value = 7727.273
$ws.Range("J51").Value = 5700
$ws.Range("L51").Value = 5700
$ws.Range("N51").Value = -6668
$ws.Range("H62").Value = 5148.9
$ws.Range("J62").Value = 12733.333
$ws.Range("L62").Value = 12733.333
$ws.Range("N62").Value = -13981.333
$ws.Range("H65").Value = 5148.9
$ws.Range("J65").Value = 12733.333
$ws.Range("L65").Value = 63666.665
$ws.Range("N65").Value = -69906.66500000001
$ws.Range("H86").Value = 4323.4
$ws.Range("I86").Value = 2005.6
$ws.Range("K86").Value = 2005.6
$ws.Range("M86").Value = -882.5999999999999
$ws.Range("H88").Value = 2001.5
$ws.Range("J88").Value = 3000
$ws.Range("L88").Value = 3000
$ws.Range("N88").Value = -3812
$ws.Range("H89").Value = 4323.4
$ws.Range("I89").Value = 2005.6
$ws.Range("K89").Value = 10028
$ws.Range("M89").Value = -4412
$ws.Range("H91").Value = 2001.5
$ws.Range("J91").Value = 3000
$ws.Range("L91").Value = 3000
$ws.Range("N91").Value = -5808
$ws.Range("H96").Value = 764693.2
$ws.Range("I96").Value = 2031.2
$ws.Range("J96").Value = 1612095.5
$ws.Range("K96").Value = 6093.6
$ws.Range("L96").Value = 4836286.5
$ws.Range("M96").Value = -4720.6
$ws.Range("N96").Value = -4839032.5
$ws.Range("H97").Value = 1740
$ws.Range("J97").Value = 1740
$ws.Range("L97").Value = 5220
$ws.Range("N97").Value = -6212
$ws.Range("H100").Value = 6392.25
$ws.Range("I100").Value = 3500.2856
$ws.Range("J100").Value = 10441
$ws.Range("K100").Value = 3500.2856
$ws.Range("L100").Value = 10441
$ws.Range("M100").Value = -2959.2856
$ws.Range("N100").Value = -11523
$ws.Range("H103").Value = 71430000
$ws.Range("I103").Value = 1000
$ws.Range("K103").Value = 3000
$ws.Range("M103").Value = -2414
$ws.Range("H107").Value = 1494.1765
$ws.Range("J107").Value = 3170.1667
$ws.Range("L107").Value = 3170.1667
$ws.Range("N107").Value = -7010.1667
$ws.Range("H116").Value = 12320.857
$ws.Range("I116").Value = 7815.3335
$ws.Range("J116").Value = 14123.066
$ws.Range("K116").Value = 7815.3335
$ws.Range("L116").Value = 14123.066
$ws.Range("M116").Value = -4373.3335
$ws.Range("N116").Value = -21007.066
$ws.Range("H132").Value = 5414.6333
$ws.Range("I132").Value = 2702.25
$ws.Range("J132").Value = 10839.4
$ws.Range("K132").Value = 8106.75
$ws.Range("L132").Value = 32518.2
$ws.Range("M132").Value = -5576.75
$ws.Range("N132").Value = -37578.2
$ws.Range("H135").Value = 1267.625
$ws.Range("I135").Value = 502.4643
$ws.Range("K135").Value = 4522.178699999999
$ws.Range("M135").Value = -1987.178699999999
$ws.Range("H137").Value = 1290.875
$ws.Range("I137").Value = 804.6667
$ws.Range("K137").Value = 2414.0001
$ws.Range("M137").Value = 135.9998999999998
$ws.Range("H138").Value = 3243.851
$ws.Range("I138").Value = 1586.5652
$ws.Range("J138").Value = 4832.0835
$ws.Range("K138").Value = 4759.6956
$ws.Range("L138").Value = 14496.2505
$ws.Range("M138").Value = 380.3044
$ws.Range("N138").Value = -24776.2505
$ws.Range("H141").Value = 33344644
$ws.Range("I141").Value = 50003756
$ws.Range("J141").Value = 26425.4
$ws.Range("K141").Value = 150011268
$ws.Range("L141").Value = 79276.20000000001
$ws.Range("M141").Value = -150006088
$ws.Range("N141").Value = -89636.20000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 157.42857
$ws.Range("J4").Value = 151
$ws.Range("L4").Value = 151
$ws.Range("N4").Value = -383
$ws.Range("H32").Value = 3996
$ws.Range("I32").Value = 4204.5454
$ws.Range("K32").Value = 4204.5454
$ws.Range("M32").Value = -3917.5454
$ws.Range("H37").Value = 50000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 50000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 50000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -50546
$ws.Range("H45").Value = 2173.8572
$ws.Range("I45").Value = 1420.6
$ws.Range("K45").Value = 1420.6
$ws.Range("M45").Value = -1043.6
$ws.Range("H61").Value = 11059461
$ws.Range("I61").Value = 15388835
$ws.Range("K61").Value = 15388835
$ws.Range("M61").Value = -15388623
$ws.Range("H63").Value = 2151.8
$ws.Range("I63").Value = 1939.75
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 1939.75
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1253.75
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2151.8
$ws.Range("I66").Value = 1939.75
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 9698.75
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -6266.75
$ws.Range("N66").Value = -21864
$ws.Range("H74").Value = 1562.4517
$ws.Range("I74").Value = 1041.5927
$ws.Range("K74").Value = 1041.5927
$ws.Range("M74").Value = -167.5926999999999
$ws.Range("H77").Value = 1562.4517
$ws.Range("I77").Value = 1041.5927
$ws.Range("K77").Value = 5207.9635
$ws.Range("M77").Value = -839.9634999999998
$ws.Range("H132").Value = 2003362.2
$ws.Range("I132").Value = 3184.8
$ws.Range("J132").Value = 20004960
$ws.Range("K132").Value = 9554.400000000001
$ws.Range("L132").Value = 60014880
$ws.Range("M132").Value = -7024.400000000001
$ws.Range("N132").Value = -60019940
$ws.Range("H136").Value = 11059461
$ws.Range("I136").Value = 15388835
$ws.Range("K136").Value = 46166505
$ws.Range("M136").Value = -46163955

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 26729.393
$ws.Range("I86").Value = 33103.24
$ws.Range("K86").Value = 33103.24
$ws.Range("M86").Value = -31980.24
$ws.Range("H89").Value = 26729.393
$ws.Range("I89").Value = 33103.24
$ws.Range("K89").Value = 165516.2
$ws.Range("M89").Value = -159900.2
$ws.Range("H94").Value = 2197.7646
$ws.Range("I94").Value = 1929.4286
$ws.Range("J94").Value = 3450
$ws.Range("K94").Value = 1929.4286
$ws.Range("L94").Value = 3450
$ws.Range("M94").Value = -1478.4286
$ws.Range("N94").Value = -4352
$ws.Range("H99").Value = 2846.8
$ws.Range("I99").Value = 2450.3333
$ws.Range("K99").Value = 2450.3333
$ws.Range("M99").Value = -952.3332999999998
$ws.Range("H107").Value = 5902.636
$ws.Range("I107").Value = 5793
$ws.Range("K107").Value = 5793
$ws.Range("M107").Value = -3873
$ws.Range("H134").Value = 7145602
$ws.Range("I134").Value = 2743.5
$ws.Range("K134").Value = 8230.5
$ws.Range("M134").Value = -5695.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 4
$ws.Range("K2").Value = 4
$ws.Range("M2").Value = 109
$ws.Range("H22").Value = 444.07144
$ws.Range("I22").Value = 406.85715
$ws.Range("K22").Value = 406.85715
$ws.Range("M22").Value = -56.85714999999999
$ws.Range("H31").Value = 38464144
$ws.Range("I31").Value = 50002550
$ws.Range("J31").Value = 2781.3333
$ws.Range("K31").Value = 50002550
$ws.Range("L31").Value = 2781.3333
$ws.Range("M31").Value = -50002255
$ws.Range("N31").Value = -3371.3333
$ws.Range("H34").Value = 38464144
$ws.Range("I34").Value = 50002550
$ws.Range("J34").Value = 2781.3333
$ws.Range("K34").Value = 50002550
$ws.Range("L34").Value = 2781.3333
$ws.Range("M34").Value = -50002348
$ws.Range("N34").Value = -3185.3333
$ws.Range("H50").Value = 42498.25
$ws.Range("J50").Value = 23331
$ws.Range("L50").Value = 23331
$ws.Range("N50").Value = -24581
$ws.Range("H51").Value = 18000
$ws.Range("I51").Value = 18000
$ws.Range("K51").Value = 18000
$ws.Range("M51").Value = -17264
$ws.Range("H58").Value = 3012.0715
$ws.Range("I58").Value = 2847.5
$ws.Range("K58").Value = 2847.5
$ws.Range("M58").Value = -2644.5
$ws.Range("H59").Value = 62499.5
$ws.Range("I59").Value = 25000
$ws.Range("K59").Value = 25000
$ws.Range("M59").Value = -23855
$ws.Range("H60").Value = 22046.5
$ws.Range("I60").Value = 22046.5
$ws.Range("K60").Value = 22046.5
$ws.Range("M60").Value = -21535.5
$ws.Range("H61").Value = 18000
$ws.Range("I61").Value = 18000
$ws.Range("K61").Value = 18000
$ws.Range("M61").Value = -17652
$ws.Range("H99").Value = 35447.727
$ws.Range("I99").Value = 11446.875
$ws.Range("J99").Value = 99450
$ws.Range("K99").Value = 11446.875
$ws.Range("L99").Value = 99450
$ws.Range("M99").Value = -9948.875
$ws.Range("N99").Value = -102446
$ws.Range("H126").Value = 35447.727
$ws.Range("I126").Value = 11446.875
$ws.Range("J126").Value = 99450
$ws.Range("K126").Value = 34340.625
$ws.Range("L126").Value = 298350
$ws.Range("M126").Value = -31870.625
$ws.Range("N126").Value = -303290
$ws.Range("H132").Value = 2812.8708
$ws.Range("I132").Value = 2718
$ws.Range("K132").Value = 8154
$ws.Range("M132").Value = -5624
$ws.Range("H134").Value = 2534.8823
$ws.Range("I134").Value = 2306.2
$ws.Range("K134").Value = 6918.599999999999
$ws.Range("M134").Value = -4383.599999999999
$ws.Range("H136").Value = 3012.0715
$ws.Range("I136").Value = 2847.5
$ws.Range("K136").Value = 8542.5
$ws.Range("M136").Value = -5992.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5598.3335
$ws.Range("I34").Value = 437.5
$ws.Range("J34").Value = 9727
$ws.Range("K34").Value = 1312.5
$ws.Range("L34").Value = 29181
$ws.Range("M34").Value = -1228.5
$ws.Range("N34").Value = -29349
$ws.Range("H39").Value = 16582.75
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 16582.75
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 49748.25
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -50336.25
$ws.Range("H55").Value = 5671.1665
$ws.Range("I55").Value = 2568.182
$ws.Range("K55").Value = 7704.545999999999
$ws.Range("M55").Value = -7527.545999999999
$ws.Range("H69").Value = 11548.333
$ws.Range("I69").Value = 656
$ws.Range("K69").Value = 1968
$ws.Range("M69").Value = -1157
$ws.Range("H72").Value = 11548.333
$ws.Range("I72").Value = 656
$ws.Range("K72").Value = 5904
$ws.Range("M72").Value = -1848
$ws.Range("H87").Value = 21666.5
$ws.Range("I87").Value = 10000
$ws.Range("K87").Value = 30000
$ws.Range("M87").Value = -28752
$ws.Range("H90").Value = 21666.5
$ws.Range("I90").Value = 10000
$ws.Range("K90").Value = 90000
$ws.Range("M90").Value = -83760
$ws.Range("H132").Value = 3027.3704
$ws.Range("I132").Value = 1562
$ws.Range("J132").Value = 4199.6665
$ws.Range("K132").Value = 14058
$ws.Range("L132").Value = 37796.9985
$ws.Range("M132").Value = -11528
$ws.Range("N132").Value = -42856.9985
$ws.Range("H137").Value = 6453.353
$ws.Range("I137").Value = 3579.4546
$ws.Range("K137").Value = 10738.3638
$ws.Range("M137").Value = -5638.363799999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8210.333000000001
$ws.Range("I3").Value = 680
$ws.Range("J3").Value = 11975.5
$ws.Range("K3").Value = 680
$ws.Range("L3").Value = 11975.5
$ws.Range("M3").Value = -564
$ws.Range("N3").Value = -12207.5
$ws.Range("H43").Value = 6895
$ws.Range("I43").Value = 6895
$ws.Range("K43").Value = 6895
$ws.Range("M43").Value = -6744
$ws.Range("H70").Value = 10022.25
$ws.Range("I70").Value = 8181.8335
$ws.Range("J70").Value = 11126.5
$ws.Range("K70").Value = 8181.8335
$ws.Range("L70").Value = 11126.5
$ws.Range("M70").Value = -7911.8335
$ws.Range("N70").Value = -11666.5
$ws.Range("H73").Value = 10022.25
$ws.Range("I73").Value = 8181.8335
$ws.Range("J73").Value = 11126.5
$ws.Range("K73").Value = 8181.8335
$ws.Range("L73").Value = 11126.5
$ws.Range("M73").Value = -7245.8335
$ws.Range("N73").Value = -12998.5
$ws.Range("H80").Value = 4215.75
$ws.Range("J80").Value = 4388.7144
$ws.Range("L80").Value = 4388.7144
$ws.Range("N80").Value = -6384.7144
$ws.Range("H83").Value = 4215.75
$ws.Range("J83").Value = 4388.7144
$ws.Range("L83").Value = 21943.572
$ws.Range("N83").Value = -31927.572
$ws.Range("H97").Value = 615.56757
$ws.Range("I97").Value = 497.0345
$ws.Range("J97").Value = 1045.25
$ws.Range("K97").Value = 497.0345
$ws.Range("L97").Value = 1045.25
$ws.Range("M97").Value = -1.03449999999998
$ws.Range("N97").Value = -2037.25
$ws.Range("H122").Value = 3281.75
$ws.Range("J122").Value = 3199.6667
$ws.Range("L122").Value = 9599.000100000001
$ws.Range("N122").Value = -14499.0001
$ws.Range("H126").Value = 1943.75
$ws.Range("I126").Value = 1900.4546
$ws.Range("K126").Value = 5701.3638
$ws.Range("M126").Value = -3231.3638
$ws.Range("H132").Value = 12502212
$ws.Range("I132").Value = 2339
$ws.Range("J132").Value = 33335334
$ws.Range("K132").Value = 7017
$ws.Range("L132").Value = 100006002
$ws.Range("M132").Value = -4487
$ws.Range("N132").Value = -100011062

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 110001
$ws.Range("I22").Value = 110001
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 110001
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -109706
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 110001
$ws.Range("I27").Value = 110001
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 110001
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -109894
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 4998
$ws.Range("I40").Value = 4998
$ws.Range("K40").Value = 4998
$ws.Range("M40").Value = -4862
$ws.Range("H55").Value = 1905.4615
$ws.Range("I55").Value = 2029.4
$ws.Range("K55").Value = 2029.4
$ws.Range("M55").Value = -1856.4
$ws.Range("H61").Value = 58827970
$ws.Range("I61").Value = 125000930
$ws.Range("J61").Value = 7555.1113
$ws.Range("K61").Value = 125000930
$ws.Range("L61").Value = 7555.1113
$ws.Range("M61").Value = -125000728
$ws.Range("N61").Value = -7959.1113
$ws.Range("H82").Value = 5274.5386
$ws.Range("I82").Value = 2742.3333
$ws.Range("K82").Value = 2742.3333
$ws.Range("M82").Value = -2381.3333
$ws.Range("H85").Value = 5274.5386
$ws.Range("I85").Value = 2742.3333
$ws.Range("K85").Value = 2742.3333
$ws.Range("M85").Value = -1494.3333
$ws.Range("H93").Value = 1545555.9
$ws.Range("I93").Value = 781.8333
$ws.Range("K93").Value = 781.8333
$ws.Range("M93").Value = 466.1667
$ws.Range("H113").Value = 58827970
$ws.Range("I113").Value = 125000930
$ws.Range("J113").Value = 7555.1113
$ws.Range("K113").Value = 125000930
$ws.Range("L113").Value = 7555.1113
$ws.Range("M113").Value = -124998760
$ws.Range("N113").Value = -11895.1113
$ws.Range("H132").Value = 2676.7058
$ws.Range("I132").Value = 1759
$ws.Range("K132").Value = 5277
$ws.Range("M132").Value = -2747
$ws.Range("H136").Value = 4035.111
$ws.Range("I136").Value = 2876.0908
$ws.Range("K136").Value = 8628.2724
$ws.Range("M136").Value = -6078.2724

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 13883.625
$ws.Range("J96").Value = 15274
$ws.Range("L96").Value = 15274
$ws.Range("N96").Value = -18020
$ws.Range("H122").Value = 2953
$ws.Range("I122").Value = 2404.1667
$ws.Range("J122").Value = 3501.8333
$ws.Range("K122").Value = 7212.500100000001
$ws.Range("L122").Value = 10505.4999
$ws.Range("M122").Value = -4762.500100000001
$ws.Range("N122").Value = -15405.4999
$ws.Range("H132").Value = 359262.53
$ws.Range("I132").Value = 2113.9524
$ws.Range("J132").Value = 1430708.2
$ws.Range("K132").Value = 6341.8572
$ws.Range("L132").Value = 4292124.6
$ws.Range("M132").Value = -3811.8572
$ws.Range("N132").Value = -4297184.6
$ws.Range("H136").Value = 272870.3
$ws.Range("I136").Value = 2337.3125
$ws.Range("K136").Value = 7011.9375
$ws.Range("M136").Value = -4461.9375
